$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format from an existing header cell (G1) so the new header
# cell reuses the same style index rather than creating a duplicate one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

$values = @(1, 1, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
